$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1698.1818
$ws.Range("I98").Value = 1880.0344
$ws.Range("J98").Value = 379.75
$ws.Range("K98").Value = 1880.0344
$ws.Range("L98").Value = 379.75
$ws.Range("M98").Value = -382.0344
$ws.Range("N98").Value = -3375.75

$ws.Range("H100").Value = 1614.2222
$ws.Range("I100").Value = 1504.1666
$ws.Range("J100").Value = 1834.3334
$ws.Range("K100").Value = 1504.1666
$ws.Range("L100").Value = 1834.3334
$ws.Range("M100").Value = -963.1666
$ws.Range("N100").Value = -2916.3334

$ws.Range("H122").Value = 1698.1818
$ws.Range("I122").Value = 1880.0344
$ws.Range("J122").Value = 379.75
$ws.Range("K122").Value = 5640.1032
$ws.Range("L122").Value = 1139.25
$ws.Range("M122").Value = -3190.1032
$ws.Range("N122").Value = -6039.25

$ws.Range("H137").Value = 30304840
$ws.Range("I137").Value = 1213.9584
$ws.Range("J137").Value = 111114510
$ws.Range("K137").Value = 3641.8752
$ws.Range("L137").Value = 333343530
$ws.Range("M137").Value = -1091.8752
$ws.Range("N137").Value = -333348630

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6187.339
$ws.Range("I32").Value = 3073.42
$ws.Range("J32").Value = 19162
$ws.Range("K32").Value = 3073.42
$ws.Range("L32").Value = 19162
$ws.Range("M32").Value = -2786.42
$ws.Range("N32").Value = -19736

$ws.Range("H61").Value = 2555.9412
$ws.Range("I61").Value = 1459.7858
$ws.Range("J61").Value = 7671.3335
$ws.Range("K61").Value = 1459.7858
$ws.Range("L61").Value = 7671.3335
$ws.Range("M61").Value = -1247.7858
$ws.Range("N61").Value = -8095.3335

$ws.Range("H74").Value = 4658.294
$ws.Range("I74").Value = 980.3077
$ws.Range("K74").Value = 980.3077
$ws.Range("M74").Value = -106.3077

$ws.Range("H77").Value = 4658.294
$ws.Range("I77").Value = 980.3077
$ws.Range("K77").Value = 4901.5385
$ws.Range("M77").Value = -533.5384999999997

$ws.Range("H122").Value = 1738.909
$ws.Range("I122").Value = 1514.2222
$ws.Range("K122").Value = 4542.6666
$ws.Range("M122").Value = -2092.6666

$ws.Range("H136").Value = 2555.9412
$ws.Range("I136").Value = 1459.7858
$ws.Range("J136").Value = 7671.3335
$ws.Range("K136").Value = 4379.357400000001
$ws.Range("L136").Value = 23014.0005
$ws.Range("M136").Value = -1829.357400000001
$ws.Range("N136").Value = -28114.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1449.2693
$ws.Range("I107").Value = 1525.1305
$ws.Range("J107").Value = 867.6667
$ws.Range("K107").Value = 1525.1305
$ws.Range("L107").Value = 867.6667
$ws.Range("M107").Value = 394.8695
$ws.Range("N107").Value = -4707.6667

$ws.Range("H134").Value = 26232.244
$ws.Range("I134").Value = 27710.404
$ws.Range("J134").Value = 5538
$ws.Range("K134").Value = 83131.212
$ws.Range("L134").Value = 16614
$ws.Range("M134").Value = -80596.212
$ws.Range("N134").Value = -21684

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1458.2727
$ws.Range("I16").Value = 1406.375
$ws.Range("K16").Value = 1406.375
$ws.Range("M16").Value = -1119.375

$ws.Range("H94").Value = 1131.1111
$ws.Range("I94").Value = 1066.6666
$ws.Range("J94").Value = 1163.3334
$ws.Range("K94").Value = 1066.6666
$ws.Range("L94").Value = 1163.3334
$ws.Range("M94").Value = -615.6666
$ws.Range("N94").Value = -2065.3334

$ws.Range("H113").Value = 1458.2727
$ws.Range("I113").Value = 1406.375
$ws.Range("K113").Value = 1406.375
$ws.Range("M113").Value = 763.625

$ws.Range("H132").Value = 2028.5135
$ws.Range("I132").Value = 1466.6
$ws.Range("J132").Value = 3199.1667
$ws.Range("K132").Value = 4399.799999999999
$ws.Range("L132").Value = 9597.500100000001
$ws.Range("M132").Value = -1869.799999999999
$ws.Range("N132").Value = -14657.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 2000
$ws.Range("J49").Value = 2000
$ws.Range("L49").Value = 6000
$ws.Range("N49").Value = -6312

$ws.Range("H68").Value = 1200
$ws.Range("I68").Value = 1166.6666
$ws.Range("J68").Value = 1225
$ws.Range("K68").Value = 3499.9998
$ws.Range("L68").Value = 3675
$ws.Range("M68").Value = -2688.9998
$ws.Range("N68").Value = -5297

$ws.Range("H71").Value = 1200
$ws.Range("I71").Value = 1166.6666
$ws.Range("J71").Value = 1225
$ws.Range("K71").Value = 10499.9994
$ws.Range("L71").Value = 11025
$ws.Range("M71").Value = -6443.999400000001
$ws.Range("N71").Value = -19137

$ws.Range("H100").Value = 1980
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H107").Value = 57501.37
$ws.Range("I107").Value = 125192.375
$ws.Range("J107").Value = 37444.777
$ws.Range("K107").Value = 375577.125
$ws.Range("L107").Value = 112334.331
$ws.Range("M107").Value = -373657.125
$ws.Range("N107").Value = -116174.331

$ws.Range("H108").Value = 308.66666
$ws.Range("I108").Value = 308.66666
$ws.Range("K108").Value = 925.9999799999999
$ws.Range("M108").Value = 1954.00002

$ws.Range("H117").Value = 3743.75
$ws.Range("I117").Value = 475
$ws.Range("J117").Value = 4833.3335
$ws.Range("K117").Value = 1425
$ws.Range("L117").Value = 14500.0005
$ws.Range("M117").Value = 2017
$ws.Range("N117").Value = -21384.0005

$ws.Range("H131").Value = 1022234.1
$ws.Range("I131").Value = 4510.75
$ws.Range("J131").Value = 1164242.1
$ws.Range("K131").Value = 13532.25
$ws.Range("L131").Value = 3492726.3
$ws.Range("M131").Value = -8492.25
$ws.Range("N131").Value = -3502806.3

$ws.Range("H132").Value = 35715230
$ws.Range("I132").Value = 43479020
$ws.Range("J132").Value = 1818
$ws.Range("K132").Value = 391311180
$ws.Range("L132").Value = 16362
$ws.Range("M132").Value = -391308650
$ws.Range("N132").Value = -21422

$ws.Range("H133").Value = 7656.722
$ws.Range("I133").Value = 4738
$ws.Range("K133").Value = 14214
$ws.Range("M133").Value = -9154

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1812.4445
$ws.Range("I126").Value = 1762.4
$ws.Range("J126").Value = 1875
$ws.Range("K126").Value = 5287.200000000001
$ws.Range("L126").Value = 5625
$ws.Range("M126").Value = -2817.200000000001
$ws.Range("N126").Value = -10565

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1388.5483
$ws.Range("I7").Value = 1164.5
$ws.Range("J7").Value = 1795.909
$ws.Range("K7").Value = 1164.5
$ws.Range("L7").Value = 1795.909
$ws.Range("M7").Value = -1052.5
$ws.Range("N7").Value = -2019.909

$ws.Range("H122").Value = 4874.523
$ws.Range("I122").Value = 5246.0347
$ws.Range("J122").Value = 4156.2666
$ws.Range("K122").Value = 15738.1041
$ws.Range("L122").Value = 12468.7998
$ws.Range("M122").Value = -13288.1041
$ws.Range("N122").Value = -17368.7998

$ws.Range("H126").Value = 1388.5483
$ws.Range("I126").Value = 1164.5
$ws.Range("J126").Value = 1795.909
$ws.Range("K126").Value = 3493.5
$ws.Range("L126").Value = 5387.727000000001
$ws.Range("M126").Value = -1023.5
$ws.Range("N126").Value = -10327.727

$ws.Range("H132").Value = 4808.057
$ws.Range("I132").Value = 5657
$ws.Range("J132").Value = 3371.3845
$ws.Range("K132").Value = 16971
$ws.Range("L132").Value = 10114.1535
$ws.Range("M132").Value = -14441
$ws.Range("N132").Value = -15174.1535

$ws.Range("H136").Value = 2298.8
$ws.Range("I136").Value = 1123.5
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 3370.5
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -820.5
$ws.Range("N136").Value = -26100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 7012.1816
$ws.Range("I136").Value = 7142
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 21426
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -18876
$ws.Range("N136").Value = -20100
